# The two Pearson logo images (in the footers) are currently named
# "image1.png" and need to become "image2.png"; the BTEC logo image
# (in the header) is currently named "image2.jpg" and needs to become
# "image1.jpg".
#
# InlineShape has no settable .Name property in the Word object model
# (only Shape/ShapeRange do), so we briefly promote each inline picture
# to a floating Shape, rename it, then convert it straight back to an
# inline shape in place.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlinePicture($range, $shapeIndex, $newName) {
    $inlineShape = $range.InlineShapes($shapeIndex)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Header: BTec_Logo-Orange -> image1.jpg
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $h = $sec.Headers($i)
    if ($h.Exists) {
        for ($j = 1; $j -le $h.Range.InlineShapes.Count; $j++) {
            $shp = $h.Range.InlineShapes($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                Rename-InlinePicture $h.Range $j "image1.jpg"
            }
        }
    }
}

# Footers: both Pearson logo pictures image1.png -> image2.png
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $f = $sec.Footers($i)
    if ($f.Exists) {
        for ($j = 1; $j -le $f.Range.InlineShapes.Count; $j++) {
            $shp = $f.Range.InlineShapes($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                Rename-InlinePicture $f.Range $j "image2.png"
            }
        }
    }
}
